# إضافة حدث جديد في Card19
#
# 1) Row 14 previously had a handful of "placeholder" cells (B..K and M)
#    that were completely blank. This backfills them with the literal
#    text "nan" - the same sentinel used everywhere else in this sheet to
#    mean "no value" - matching the convention used by every other row.
# 2) A brand-new service-event row (row 15) is appended for card 19 with
#    a date, a correction note and who serviced it; the Min/Max Tones and
#    wear-indicator columns (B..K, M) are left blank for this event, same
#    as row 14 above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

$blankCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "M")

# --- Row 14: backfill the empty placeholder cells with "nan" ---------------
foreach ($col in $blankCols) {
    $ws.Range($col + "14").Value = "nan"
}

# --- Row 15: new service event for card 19 ----------------------------------
# Column A holds the card number "19" as text (matching every other row in
# this column). A leading apostrophe forces text interpretation so "19"
# isn't auto-converted to a number; likewise a lone apostrophe leaves the
# Min/Max-Tones/indicator cells as blank text, matching the source sheet's
# convention of never leaving a truly empty cell behind.
$ws.Range("A15").Value = "'19"
foreach ($col in $blankCols) {
    $ws.Range($col + "15").Value = "'"
}
$ws.Range("L15").Value = "5\12\2024"
$ws.Range("N15").Value = "تم سن الفلاتس"
$ws.Range("O15").Value = "الخبير"

# Row 14 is plain, unstyled text; copy its formatting onto row 15 so the
# new cells don't keep the "text quote-prefix" styling that the leading
# apostrophes above would otherwise leave behind - row 15 ends up
# formatted exactly like the row above it.
$ws.Range("A14:O14").Copy()
$ws.Range("A15:O15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
